$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "2025-10-09 01:24:12"
$ws.Range("B3").Value = "Minsk"
$ws.Range("C3").Value = 6.88
$ws.Range("D3").Value = 5.81
$ws.Range("E3").Value = "пасмурно"
$ws.Range("F3").Value = 93
$ws.Range("G3").Value = 1.75
